$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A slightly (13.85546875 -> 14.85546875) to fit the longer
# date string in the new row ("8/19/16 18:35" vs "8/5/16 19:06")
$ws.Columns.Item(1).ColumnWidth = 14

# Add new row 3 of data
$ws.Range("A3").Value = Get-Date -Year 2016 -Month 8 -Day 19 -Hour 18 -Minute 35 -Second 35
$ws.Range("B3").Value = "Gilead Sciences, Inc."
$ws.Range("C3").Value = "GILD"
$ws.Range("D3").Value = 80.91
$ws.Range("E3").Value = 80.7
$ws.Range("F3").Value = -0.14000000000000001
